# Apply updated cryptocurrency price/volume data to Sheet1
# (values are stored as text in the source data, e.g. "73.73" or "41.954.73",
#  so numeric-looking Price values are written via a Text number format to
#  avoid Excel auto-converting them to numbers)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.061.64"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.230.51"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.62"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.73"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -4.55%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.60"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.37%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.95"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -5.39%  "
$ws.Range("D14").Value = "2.564.72"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.34"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").Value = "2.222.84"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "41.882.39"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000106"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.25"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.90"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.23"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +12.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.64"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -7.95%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.42"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.62"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.12"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.56"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.63"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0801"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.06"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -9.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -7.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.46"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -7.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.06"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.67"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.198"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.76"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.97"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "2.432.42"
$ws.Range("E51").Value = "  -1.46%  "
